# Update specific numeric values in Sheet1 ("Update Name of Algo")
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "A3" = -21.933
    "C6" = -11.759
    "A14" = -21.727
    "C18" = -12.484
    "C19" = -12.363
    "A21" = -20.682
    "B22" = 7.131
    "A23" = -20.766
    "B24" = 5.454000000000001
    "A25" = -21.571
    "A26" = -21.607
    "B28" = 5.83
    "A29" = -21.577
    "B36" = 6.965999999999999
    "C44" = -12.74
    "B45" = 5.914
    "C47" = -12.382
    "B48" = 5.48
    "B49" = 6.237
    "C51" = -11.589
    "B52" = 4.992
    "A53" = -20.293
    "B53" = 7.512
    "B54" = 5.702
    "A57" = -22.15
    "C57" = -13.667
    "A59" = -22.461
    "C64" = -10.76
    "A69" = -21.591
    "B70" = 4.935
    "A79" = -21.2
    "C80" = -12.55
    "A83" = -21.884
    "B86" = 5.583
    "B87" = 5.336
    "B89" = 5.95
    "A91" = -20.987
    "C92" = -10.975
    "A93" = -21.687
    "C94" = -11.35
    "C96" = -10.832
    "B101" = 5.116999999999999
    "C101" = -13.138
    "A103" = -21.948
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
